$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header row formatting down into the two new rows first so that
# they reuse the existing cell style (s="2") instead of creating duplicates.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$timeTypes = "dateTime" + [char]0x0135 + ", Period" + [char]0x0135 + ", Timing" + [char]0x0135 + ", instant" + [char]0x0135
$valueTypes = "string" + [char]0x0135 + ", CodeableConcept" + [char]0x0135 + ", Quantity" + [char]0x0135 + ", boolean" + [char]0x0135 + ", integer" + [char]0x0135 + ", Range" + [char]0x0135 + ", Ratio" + [char]0x0135 + ", SampledData" + [char]0x0135 + ", time" + [char]0x0135 + ", dateTime" + [char]0x0135 + ", Period" + [char]0x0135

$row2 = @(
    "us-core-treatment-intervention-preference-bindings",
    "US Core Treatment Intervention Preference Bindings Profile",
    "null#treatment-intervention-preference",
    "",
    "LOINC#75773-2",
    "",
    $timeTypes,
    $valueTypes,
    "optional",
    "",
    ""
)

$row3 = @(
    "us-core-treatment-intervention-preference-grouping",
    "US Core Treatment Intervention Preference Grouping Profile",
    "null#treatment-intervention-preference",
    "",
    "LOINC#75773-2",
    "",
    $timeTypes,
    $valueTypes,
    "optional",
    "",
    ""
)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}
